# Updates "想去人数" (F) and "最低票价" (G) figures across all sheets to match
# the refreshed scrape output (gh-pages data regeneration at commit 456a3b4).

$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# ---- Sheet: 展览 (Exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
Set-Cell $ws 2  6 8303   # F2
Set-Cell $ws 2  7 95     # G2
Set-Cell $ws 3  6 131    # F3
Set-Cell $ws 4  6 103    # F4
Set-Cell $ws 4  7 70     # G4
Set-Cell $ws 5  6 35737  # F5
Set-Cell $ws 7  6 612    # F7
Set-Cell $ws 8  6 727    # F8
Set-Cell $ws 12 6 816    # F12
Set-Cell $ws 15 6 448    # F15
Set-Cell $ws 17 6 585    # F17
Set-Cell $ws 18 6 162    # F18
Set-Cell $ws 19 6 433    # F19
Set-Cell $ws 20 6 428    # F20
Set-Cell $ws 21 6 1127   # F21
Set-Cell $ws 23 6 745    # F23
Set-Cell $ws 24 6 2400   # F24
Set-Cell $ws 25 6 886    # F25
Set-Cell $ws 26 6 510    # F26
Set-Cell $ws 27 6 79     # F27
Set-Cell $ws 30 6 673    # F30
Set-Cell $ws 31 6 673    # F31
Set-Cell $ws 32 6 15     # F32
Set-Cell $ws 33 6 1103   # F33

# ---- Sheet: 演出 (Performances) ----
$ws = $wb.Worksheets.Item("演出")
Set-Cell $ws 2  6 302    # F2
Set-Cell $ws 10 6 3      # F10

# ---- Sheet: 本地生活 (Local Life) ----
$ws = $wb.Worksheets.Item("本地生活")
Set-Cell $ws 2 6 561     # F2

# ---- Sheet: 全部类型 (All Types, combined listing) ----
$ws = $wb.Worksheets.Item("全部类型")
Set-Cell $ws 2  6 561    # F2
Set-Cell $ws 3  6 8303   # F3
Set-Cell $ws 3  7 95     # G3
Set-Cell $ws 4  6 131    # F4
Set-Cell $ws 5  6 103    # F5
Set-Cell $ws 5  7 70     # G5
Set-Cell $ws 6  6 302    # F6
Set-Cell $ws 7  6 35737  # F7
Set-Cell $ws 9  6 612    # F9
Set-Cell $ws 10 6 727    # F10
Set-Cell $ws 18 6 816    # F18
Set-Cell $ws 21 6 448    # F21
Set-Cell $ws 26 6 3      # F26
Set-Cell $ws 28 6 585    # F28
Set-Cell $ws 29 6 162    # F29
Set-Cell $ws 30 6 433    # F30
Set-Cell $ws 31 6 428    # F31
Set-Cell $ws 32 6 1127   # F32
Set-Cell $ws 34 6 745    # F34
Set-Cell $ws 35 6 2400   # F35
Set-Cell $ws 36 6 886    # F36
Set-Cell $ws 37 6 510    # F37
Set-Cell $ws 38 6 79     # F38
Set-Cell $ws 42 6 673    # F42
Set-Cell $ws 43 6 673    # F43
Set-Cell $ws 44 6 15     # F44
Set-Cell $ws 45 6 1103   # F45
